$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell X1
$ws.Range("X1").Value = "Unnamed: 23"
$ws.Range("X1").Style = $ws.Range("W1").Style

# New row 3 data
$ws.Range("A3").Value = "09/12/2025"
$ws.Range("B3").Value = "Qwen2.5-14B-Instruct"
$ws.Range("C3").Value = 0.4739130434782609
$ws.Range("D3").Value = 0.367003367003367
$ws.Range("E3").Value = 0.4136622390891841
$ws.Range("F3").Value = 0.5725120977878331
$ws.Range("G3").Value = 0.3548946330167079
$ws.Range("H3").Value = 0.4149033099591284
$ws.Range("I3").Value = 0.5673708553610514
$ws.Range("J3").Value = 0.367003367003367
$ws.Range("K3").Value = 0.4269194440240616
$ws.Range("L3").Value = 0.5391304347826087
$ws.Range("M3").Value = 0.4189189189189189
$ws.Range("N3").Value = 0.4714828897338403
$ws.Range("O3").Value = 124
$ws.Range("P3").Value = 106
$ws.Range("Q3").Value = 172
$ws.Range("R3").Value = 297
$ws.Range("S3").Value = 0.9487118034751348
$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-14B-Instruct_2_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-14B-Instruct_2_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.024 kWh"
$ws.Range("X3").Value = ""
